$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue #13: allow two columns in a metadata file to be related to each
# other so hierarchical SKOS concepts can be built. The metadata sheet
# gains a "slug" row (machine-readable column id, no namespace prefix) in
# row 2, the old namespaced ids shift down to row 3, the existing "medida"
# row becomes row 4, and a brand-new "xsd:*" datatype row is appended as
# row 5.

# Row 2 — short slug identifying each column (used to relate columns).
$ws.Range("A2").Value = "territorio"
$ws.Range("B2").Value = "case--when-numero-trabajadores-empr--10-then-01-de-1-a-9-afiliados-when-numero-trabajadores-empr--20-then-02-de-10-a-19-afiliados-when-numero-trabajadores-empr--50-then-03-de-20-a-49-afiliados-when-empre"
$ws.Range("C2").Value = "estrato"
$ws.Range("D2").Value = "numero-empresas"
$ws.Range("E2").Value = "direccion-provincial-nombre"
$ws.Range("F2").Value = "mes-y-ano"
$ws.Range("G2").Value = "direccion-provincial-codigo"

# Row 3 — the namespaced "iaest-measure:" id that used to live in row 2.
$ws.Range("A3").Value = "iaest-measure:territorio"
$ws.Range("B3").Value = "iaest-measure:case--when-numero-trabajadores-empr--10-then-01-de-1-a-9-afiliados-when-numero-trabajadores-empr--20-then-02-de-10-a-19-afiliados-when-numero-trabajadores-empr--50-then-03-de-20-a-49-afiliados-when-empre"
$ws.Range("C3").Value = "iaest-measure:estrato"
$ws.Range("D3").Value = "iaest-measure:numero-empresas"
$ws.Range("E3").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("F3").Value = "iaest-measure:mes-y-ano"
$ws.Range("G3").Value = "null"

# Row 4 — every column is flagged as a "medida" (measure), except the
# last column, which stays "null".
$ws.Range("A4").Value = "medida"
$ws.Range("B4").Value = "medida"
$ws.Range("C4").Value = "medida"
$ws.Range("D4").Value = "medida"
$ws.Range("E4").Value = "medida"
$ws.Range("F4").Value = "medida"
$ws.Range("G4").Value = "null"

# Row 5 (new) — XSD datatype per column.
$ws.Range("A5").Value = "xsd:string"
$ws.Range("B5").Value = "xsd:string"
$ws.Range("C5").Value = "xsd:string"
$ws.Range("D5").Value = "xsd:int"
$ws.Range("E5").Value = "xsd:string"
$ws.Range("F5").Value = "xsd:string"
$ws.Range("G5").Value = "null"

# Give the new row the same formatting as the rest of the data rows.
$ws.Range("A4:G4").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)
